# Add a new "2022" year column (P) to the annual data table on the sheet,
# matching the formatting of the existing 2021 column (O), and refresh a
# few of the most-recent data-row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (thin separator/border row above the header) ---
# O3 just carries a bottom border with no value; replicate that border on P3.
$ws.Range("P3").Borders.Item(9).Weight = -4138

# --- Row 4 (year headers) ---
# Copy O4's format (bordered/bold year header) onto P4, then set its value.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# --- Row 5 (data values) ---
# Refresh the last three existing years and add the new 2022 value.
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 2.6

$ws.Application.CutCopyMode = $false

# Move/restore the active selection the way the author left it.
$ws.Range("P3").Select()
